$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Original reference list")

# Insert a new row at position 26, shifting existing rows 26+ down
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = "Ferrer-Cervantes, Mendez-Gonzalez, Quintana-Ascencio, Dorantes, Dzib & Duran"
$ws.Range("B26").Value = "Popul Ecol"
$ws.Range("C26").Value = 2012
$ws.Range("D26").Value = "Planta"
$ws.Range("E26").Value = "Mammillaria gaumeri"
$ws.Range("G26").Value = "Ferrer-Cervantes, M.E., Méndez-González, M.E., Quintana-Ascencio, P.-F., Dorantes, A., Dzib, G. & Durán, R. (2012) Population dynamics of the cactus Mammillaria gaumeri: an integral projection model approach. Population Ecology, 54, 321-334."

$ws.Range("E26").Font.Italic = $true
